# Update the "Out of PO" roster sheet:
#  - remove the obsolete 19th row (table shrinks from 18 to 17 data rows)
#  - rewrite the player / position / team table with the refreshed roster

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The new roster only has 17 data rows (rows 2-18), one fewer than before.
# Delete the now-unused last row first so the dimension shrinks correctly.
$ws.Range("A19:C19").Delete()

$ws.Range("A2").Value = "Kelly Oubre Jr."
$ws.Range("B2").Value = "SG,SF"
$ws.Range("C2").Value = "Philadelphia 76ers"

$ws.Range("A3").Value = "Malik Beasley"
$ws.Range("B3").Value = "SG,SF"
$ws.Range("C3").Value = "Detroit Pistons"

$ws.Range("A4").Value = "Josh Hart"
$ws.Range("B4").Value = "SG,SF,PF"
$ws.Range("C4").Value = "New York Knicks"

$ws.Range("A5").Value = "De'Andre Hunter"
$ws.Range("B5").Value = "SF,PF"
$ws.Range("C5").Value = "Cleveland Cavaliers"

$ws.Range("A6").Value = "Andrew Wiggins"
$ws.Range("B6").Value = "SF,PF"
$ws.Range("C6").Value = "Miami Heat"

$ws.Range("A7").Value = "Victor Wembanyama"
$ws.Range("B7").Value = "C"
$ws.Range("C7").Value = "San Antonio Spurs"

$ws.Range("A8").Value = "Domantas Sabonis"
$ws.Range("B8").Value = "C"
$ws.Range("C8").Value = "Sacramento Kings"

$ws.Range("A9").Value = "Kel'el Ware"
$ws.Range("B9").Value = "PF,C"
$ws.Range("C9").Value = "Miami Heat"

$ws.Range("A10").Value = "Kristaps Porzingis"
$ws.Range("B10").Value = "PF,C"
$ws.Range("C10").Value = "Boston Celtics"

$ws.Range("A11").Value = "Donovan Mitchell"
$ws.Range("B11").Value = "PG,SG"
$ws.Range("C11").Value = "Cleveland Cavaliers"

$ws.Range("A12").Value = "Jaden McDaniels"
$ws.Range("B12").Value = "SF,PF"
$ws.Range("C12").Value = "Minnesota Timberwolves"

$ws.Range("A13").Value = "Michael Porter Jr."
$ws.Range("B13").Value = "SF,PF"
$ws.Range("C13").Value = "Denver Nuggets"

$ws.Range("A14").Value = "Naji Marshall"
$ws.Range("B14").Value = "SG,SF"
$ws.Range("C14").Value = "Dallas Mavericks"

$ws.Range("A15").Value = "Dyson Daniels"
$ws.Range("B15").Value = "PG,SG,SF"
$ws.Range("C15").Value = "Atlanta Hawks"

$ws.Range("A16").Value = "Bilal Coulibaly"
$ws.Range("B16").Value = "SG,SF"
$ws.Range("C16").Value = "Washington Wizards"

$ws.Range("A17").Value = "Cam Thomas"
$ws.Range("B17").Value = "SG,SF"
$ws.Range("C17").Value = "Brooklyn Nets"

$ws.Range("A18").Value = "Alperen Sengün"
$ws.Range("B18").Value = "C"
$ws.Range("C18").Value = "Houston Rockets"
